# 201708044226刘研.xlsx - add a new daily-log entry (刘广部,胡颖慧,崔梦婷 DaoTest
# AND ServiceTest均通过) dated 2019-5-17 10:07:09, pushing the two previous
# "Service层ManagerTestCase的建立 / 00:00--00:30" values down into their own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 keeps the text that used to live there before the insert
# ("Service层ManagerTestCase的建立" / "00:00--00:30") - row 55 (李博文...,
# 20:30--21:00) is untouched.
$ws.Range("C56").Value = "Service层ManagerTestCase的建立"
$ws.Range("D56").Value = "00:00--00:30"

# Brand new row 57 with the new log entry.
$ws.Range("A57").Value = "2019年5月17日10:04:46"
$ws.Range("B57").Value = "周五"
$ws.Range("C57").Value = "刘广部，胡颖慧，崔梦婷DaoTest AND ServiceTest均通过"

# The time-range cell picked up a time number format (h:mm, i.e. numFmtId 20)
# while still holding literal text, matching the source workbook.
$ws.Range("D57").NumberFormat = "h:mm"
$ws.Range("D57").Value = "09:00--10:00"

# Leave the selection where the author ended up after typing the new row.
$ws.Range("C58").Select()
